{"js": "const replacements = [\n    [\"57\u00d733=\", \"83\u00d747=\"],\n    [\"29\u00d780=\", \"79\u00d738=\"],\n    [\"39\u00d762=\", \"70\u00d777=\"],\n    [\"41\u00d725=\", \"76\u00d774=\"],\n    [\"58\u00d779=\", \"11\u00d749=\"],\n    [\"38\u00d748=\", \"23\u00d742=\"],\n    [\"39\u00d731=\", \"65\u00d719=\"],\n    [\"88\u00d713=\", \"27\u00d758=\"],\n    [\"40\u00d720=\", \"75\u00d756=\"],\n    [\"94\u00d714=\", \"92\u00d776=\"],\n    [\"11\u00d760=\", \"53\u00d792=\"],\n    [\"53\u00d755=\", \"75\u00d735=\"],\n    [\"55\u00d712=\", \"44\u00d726=\"],\n    [\"21\u00d748=\", \"35\u00d730=\"],\n    [\"68\u00d712=\", \"18\u00d774=\"],\n    [\"86\u00d753=\", \"88\u00d781=\"],\n    [\"59\u00d789=\", \"98\u00d779=\"],\n    [\"11\u00d732=\", \"23\u00d749=\"],\n    [\"96\u00d797=\", \"96\u00d741=\"],\n    [\"29\u00d791=\", \"67\u00d722=\"],\n    [\"90\u00d759=\", \"91\u00d762=\"],\n    [\"50\u00d789=\", \"16\u00d753=\"],\n    [\"77\u00d788=\", \"62\u00d754=\"],\n    [\"50\u00d730=\", \"30\u00d740=\"],\n    [\"32\u00d759=\", \"16\u00d718=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"57\u00d733=\"; New = \"83\u00d747=\" },\n    @{ Old = \"29\u00d780=\"; New = \"79\u00d738=\" },\n    @{ Old = \"39\u00d762=\"; New = \"70\u00d777=\" },\n    @{ Old = \"41\u00d725=\"; New = \"76\u00d774=\" },\n    @{ Old = \"58\u00d779=\"; New = \"11\u00d749=\" },\n    @{ Old = \"38\u00d748=\"; New = \"23\u00d742=\" },\n    @{ Old = \"39\u00d731=\"; New = \"65\u00d719=\" },\n    @{ Old = \"88\u00d713=\"; New = \"27\u00d758=\" },\n    @{ Old = \"40\u00d720=\"; New = \"75\u00d756=\" },\n    @{ Old = \"94\u00d714=\"; New = \"92\u00d776=\" },\n    @{ Old = \"11\u00d760=\"; New = \"53\u00d792=\" },\n    @{ Old = \"53\u00d755=\"; New = \"75\u00d735=\" },\n    @{ Old = \"55\u00d712=\"; New = \"44\u00d726=\" },\n    @{ Old = \"21\u00d748=\"; New = \"35\u00d730=\" },\n    @{ Old = \"68\u00d712=\"; New = \"18\u00d774=\" },\n    @{ Old = \"86\u00d753=\"; New = \"88\u00d781=\" },\n    @{ Old = \"59\u00d789=\"; New = \"98\u00d779=\" },\n    @{ Old = \"11\u00d732=\"; New = \"23\u00d749=\" },\n    @{ Old = \"96\u00d797=\"; New = \"96\u00d741=\" },\n    @{ Old = \"29\u00d791=\"; New = \"67\u00d722=\" },\n    @{ Old = \"90\u00d759=\"; New = \"91\u00d762=\" },\n    @{ Old = \"50\u00d789=\"; New = \"16\u00d753=\" },\n    @{ Old = \"77\u00d788=\"; New = \"62\u00d754=\" },\n    @{ Old = \"50\u00d730=\"; New = \"30\u00d740=\" },\n    @{ Old = \"32\u00d759=\"; New = \"16\u00d718=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute([ref]$pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$pair.New, 2)\n}\n"}
